# edit.ps1 - reproduce the "reopened the deck" commit:
#   - the auto-updating date placeholders (datetimeFigureOut fields) on the
#     Slide Master and every Slide Layout refresh from 31/07/2024 to 20/10/2024
#     (the Notes Master's date field is untouched, matching the source diff)
#   - slide 5's title run gets merged/retyped: "Distributed ML – DL & Ray Cluster"
#     -> "Distributed ML – DL & RAY Cluster" (single run, dirty="0")
#   - slide 10's title gets a hyphen swapped for an en dash:
#     "DATA PIPELINE - NEO4J" -> "DATA PIPELINE – NEO4J"

$p = $ppt.ActivePresentation
$enDash = [char]0x2013

# ---------------------------------------------------------------------------
# 1) Refresh the "31/07/2024" date placeholders -> "20/10/2024" on the
#    slide master and every custom (slide) layout.
# ---------------------------------------------------------------------------
function Update-DateShapes($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $sh = $shapes.Item($j)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "31/07/2024") {
                $tr.Text = "20/10/2024"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShapes $master.Shapes

$layouts = $master.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DateShapes $layouts.Item($i).Shapes
}

# ---------------------------------------------------------------------------
# 2) Slide 5 title: "Distributed ML – DL & Ray Cluster" -> "... RAY Cluster"
#    merged back into a single run. The shape currently holds two runs
#    ("Distributed ML – DL & Ray " / "Cluster"); first collapse the range to
#    match the second run exactly so the formatting that survives the merge
#    is the one carrying dirty="0", then apply the real final text.
# ---------------------------------------------------------------------------
$slide5 = $p.Slides.Item(5)
$titleShape5 = $slide5.Shapes.Item(1)
$tr5 = $titleShape5.TextFrame.TextRange
$tr5.Text = "Cluster"
$tr5.Text = "Distributed ML " + $enDash + " DL & RAY Cluster"

# ---------------------------------------------------------------------------
# 3) Slide 10 title: "DATA PIPELINE - NEO4J" -> "DATA PIPELINE – NEO4J"
# ---------------------------------------------------------------------------
$slide10 = $p.Slides.Item(10)
$titleShape10 = $slide10.Shapes.Item(1)
$titleShape10.TextFrame.TextRange.Text = "DATA PIPELINE " + $enDash + " NEO4J"

Write-Output "edit.ps1 applied"
